$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("O2").Value = "2024-04-05"
$ws1.Range("Q2").Value = "2024-04-05 06:58:56 PM"
$ws1.Range("AD2").Value = "2024-04-05"
$ws1.Range("BB2").Value = "CT: Fri, Apr 05, 2024 at 7:01 PM"

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("O2").Value = "2024-04-05"
$ws2.Range("Q2").Value = "2024-04-05 07:06:33 PM"
$ws2.Range("AD2").Value = "2024-04-05"
$ws2.Range("BB2").Value = "CT: Fri, Apr 05, 2024 at 7:09 PM"

$ws5 = $wb.Worksheets.Item("Sheet5")
$ws5.Range("BB2").Value = "CT: Fri, Apr 05, 2024 at 7:14 PM"
